$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.016554117202759
$ws.Range("B1").Value = 1.721224308013916
$ws.Range("C1").Value = 4.565145969390869
$ws.Range("D1").Value = 5.298650741577148
$ws.Range("E1").Value = 1.618422508239746
